$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Mac-Address rows (regcntr_id, machine_id, device_id) to append after
# the existing data (rows 2..146); columns D-H repeat the constant values
# used throughout the table (eng / TRUE / superadmin / now() / now()).
$newRows = @(
    @(10001, 10030, 3000166),
    @(10001, 10030, 3000167),
    @(10001, 10030, 3000168),
    @(10001, 10030, 3000169),
    @(10001, 10030, 3000170),
    @(10001, 10031, 3000171),
    @(10001, 10031, 3000172),
    @(10001, 10031, 3000173),
    @(10001, 10031, 3000174),
    @(10001, 10031, 3000175)
)

$startRow = 147
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Match the saved view state: selection on the cell right after the new
# data, scrolled so the last rows are visible.
$ws.Range("A148").Select()
